$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.176022
$ws.Range("H2").Value = 0.5280659999999999
$ws.Range("I2").Value = 0.03293066697281707
$ws.Range("J2").Value = 0.03293066697281707
$ws.Range("M2").Value = 0.484733
$ws.Range("N2").Value = 1.454199
$ws.Range("O2").Value = 0.00792098608860474
$ws.Range("P2").Value = 0.00792098608860474
$ws.Range("Q2").Value = 0.08532367212599999
$ws.Range("R2").Value = 0.7679130491339999
$ws.Range("S2").Value = 0.0002608433549801596
$ws.Range("T2").Value = 0.0002608433549801595
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.176022
$ws.Range("H3").Value = 0.5280659999999999
$ws.Range("I3").Value = 0.03293066697281707
$ws.Range("J3").Value = 0.03293066697281707
$ws.Range("O3").Value = 0.1147190689515559
$ws.Range("P3").Value = 0.1147190689515559
$ws.Range("Q3").Value = 1.235736575766
$ws.Range("R3").Value = 11.121629181894
$ws.Range("S3").Value = 0.003777775455075327
$ws.Range("T3").Value = 0.003777775455075326
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.176022
$ws.Range("H4").Value = 0.5280659999999999
$ws.Range("I4").Value = 0.03293066697281707
$ws.Range("J4").Value = 0.03293066697281707
$ws.Range("M4").Value = 53.289524
$ws.Range("N4").Value = 159.868572
$ws.Range("O4").Value = 0.8708001689019901
$ws.Range("P4").Value = 0.8708001689019901
$ws.Range("Q4").Value = 9.380128593527999
$ws.Range("R4").Value = 84.42115734175199
$ws.Range("S4").Value = 0.02867603036198429
$ws.Range("T4").Value = 0.02867603036198429
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.176022
$ws.Range("H5").Value = 0.5280659999999999
$ws.Range("I5").Value = 0.03293066697281707
$ws.Range("J5").Value = 0.03293066697281707
$ws.Range("M5").Value = 0.4014323333333333
$ws.Range("N5").Value = 1.204297
$ws.Range("O5").Value = 0.006559776057849319
$ws.Range("P5").Value = 0.006559776057849319
$ws.Range("Q5").Value = 0.07066092217799999
$ws.Range("R5").Value = 0.6359482996019998
$ws.Range("S5").Value = 0.0002160178007772947
$ws.Range("T5").Value = 0.0002160178007772947
$ws.Range("I6").Value = 0.8002039325901205
$ws.Range("J6").Value = 0.8002039325901203
$ws.Range("M6").Value = 0.484733
$ws.Range("N6").Value = 1.454199
$ws.Range("O6").Value = 0.00792098608860474
$ws.Range("P6").Value = 0.00792098608860474
$ws.Range("Q6").Value = 2.073336019419667
$ws.Range("R6").Value = 18.660024174777
$ws.Range("S6").Value = 0.006338404218093149
$ws.Range("T6").Value = 0.006338404218093148
$ws.Range("I7").Value = 0.8002039325901205
$ws.Range("J7").Value = 0.8002039325901203
$ws.Range("O7").Value = 0.1147190689515559
$ws.Range("P7").Value = 0.1147190689515559
$ws.Range("S7").Value = 0.09179865011811221
$ws.Range("T7").Value = 0.09179865011811221
$ws.Range("I8").Value = 0.8002039325901205
$ws.Range("J8").Value = 0.8002039325901203
$ws.Range("M8").Value = 53.289524
$ws.Range("N8").Value = 159.868572
$ws.Range("O8").Value = 0.8708001689019901
$ws.Range("P8").Value = 0.8708001689019901
$ws.Range("Q8").Value = 227.9339132407507
$ws.Range("R8").Value = 2051.405219166756
$ws.Range("S8").Value = 0.6968177196555135
$ws.Range("T8").Value = 0.6968177196555135
$ws.Range("I9").Value = 0.8002039325901205
$ws.Range("J9").Value = 0.8002039325901203
$ws.Range("M9").Value = 0.4014323333333333
$ws.Range("N9").Value = 1.204297
$ws.Range("O9").Value = 0.006559776057849319
$ws.Range("P9").Value = 0.006559776057849319
$ws.Range("Q9").Value = 1.717036215936778
$ws.Range("R9").Value = 15.453325943431
$ws.Range("S9").Value = 0.005249158598401542
$ws.Range("T9").Value = 0.005249158598401542
$ws.Range("G10").Value = 0.891934
$ws.Range("H10").Value = 2.675802
$ws.Range("I10").Value = 0.1668654004370625
$ws.Range("J10").Value = 0.1668654004370625
$ws.Range("M10").Value = 0.484733
$ws.Range("N10").Value = 1.454199
$ws.Range("O10").Value = 0.00792098608860474
$ws.Range("P10").Value = 0.00792098608860474
$ws.Range("Q10").Value = 0.432349843622
$ws.Range("R10").Value = 3.891148592598
$ws.Range("S10").Value = 0.001321738515531432
$ws.Range("T10").Value = 0.001321738515531432
$ws.Range("G11").Value = 0.891934
$ws.Range("H11").Value = 2.675802
$ws.Range("I11").Value = 0.1668654004370625
$ws.Range("J11").Value = 0.1668654004370625
$ws.Range("O11").Value = 0.1147190689515559
$ws.Range("P11").Value = 0.1147190689515559
$ws.Range("Q11").Value = 6.261691532702
$ws.Range("R11").Value = 56.355223794318
$ws.Range("S11").Value = 0.01914264337836837
$ws.Range("T11").Value = 0.01914264337836837
$ws.Range("G12").Value = 0.891934
$ws.Range("H12").Value = 2.675802
$ws.Range("I12").Value = 0.1668654004370625
$ws.Range("J12").Value = 0.1668654004370625
$ws.Range("M12").Value = 53.289524
$ws.Range("N12").Value = 159.868572
$ws.Range("O12").Value = 0.8708001689019901
$ws.Range("P12").Value = 0.8708001689019901
$ws.Range("Q12").Value = 47.530738299416
$ws.Range("R12").Value = 427.776644694744
$ws.Range("S12").Value = 0.1453064188844923
$ws.Range("T12").Value = 0.1453064188844923
$ws.Range("G13").Value = 0.891934
$ws.Range("H13").Value = 2.675802
$ws.Range("I13").Value = 0.1668654004370625
$ws.Range("J13").Value = 0.1668654004370625
$ws.Range("M13").Value = 0.4014323333333333
$ws.Range("N13").Value = 1.204297
$ws.Range("O13").Value = 0.006559776057849319
$ws.Range("P13").Value = 0.006559776057849319
$ws.Range("Q13").Value = 0.3580511467993334
$ws.Range("R13").Value = 3.222460321194
$ws.Range("S13").Value = 0.001094599658670482
$ws.Range("T13").Value = 0.001094599658670482
